{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map of 0-based row index -> new cell text (single-column table).\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"1675\",\n  5: \"0.00065\",\n  6: \"0.00017\",\n  7: \"0.00005\",\n  8: \"0.00026\",\n  9: \"0.00030\",\n  10: \"0.00042\",\n  11: \"0.31291\",\n  // These rows previously held a full tab-separated breakdown line;\n  // collapse each back down to just its first (summary) value.\n  43: \"99.85\",\n  44: \"0.31\",\n  45: \"203\",\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const cell = table.getCell(Number(rowIndex), 0);\n  cell.value = updates[rowIndex];\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfunction Set-CellText($rowIndex, $text) {\n    $cell = $t.Rows.Item($rowIndex).Cells.Item(1)\n    $cell.Range.Text = $text\n}\n\n# Summary-statistics rows near the top of the table (1-based row numbers)\nSet-CellText 1 \"0M\"\nSet-CellText 2 \"0M\"\nSet-CellText 3 \"0M\"\nSet-CellText 4 \"1675\"\nSet-CellText 6 \"0.00065\"\nSet-CellText 7 \"0.00017\"\nSet-CellText 8 \"0.00005\"\nSet-CellText 9 \"0.00026\"\nSet-CellText 10 \"0.00030\"\nSet-CellText 11 \"0.00042\"\nSet-CellText 12 \"0.31291\"\n\n# The last three rows previously held a full tab-separated breakdown line;\n# collapse each back down to just its first (summary) value.\nSet-CellText 44 \"99.85\"\nSet-CellText 45 \"0.31\"\nSet-CellText 46 \"203\"\n"}
